$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.310.21"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.930.49"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7572"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.07"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.81%  "
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07035"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7800"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08024"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "1.940.76"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.399"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "30.320.28"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "253.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007964"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D21").Value = "2.190.73"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.719"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.520"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("E28").Value = "  +4.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.213"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.370"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.519"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.411"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.145"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05232"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.316"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7545"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.788"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01952"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.805"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.498"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4494"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.971"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8358"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.957"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.592"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "980.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1207"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.35%  "
